# "Fases Digestion tamaño de particula"
# - D column ("tiempo") switches from raw minute counts to fractions of an
#   hour, expressed as formulas "=<minutes>/60" (D6 is left untouched, and
#   D13 becomes a plain literal 4 instead of a formula).
# - The sheet view no longer freezes/scrolls to topLeftCell B1 and the
#   active selection moves from H2:H13 to E26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula  = "=0/60"
$ws.Range("D3").Formula  = "=10/60"
$ws.Range("D4").Formula  = "=60/60"
$ws.Range("D5").Formula  = "=120/60"
$ws.Range("D7").Formula  = "=150/60"
$ws.Range("D8").Formula  = "=165/60"
$ws.Range("D9").Formula  = "=180/60"
$ws.Range("D10").Formula = "=195/60"
$ws.Range("D11").Formula = "=210/60"
$ws.Range("D12").Formula = "=225/60"
$ws.Range("D13").Value   = 4

[void]$ws.Range("E26").Select()
